$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = "'0.57"
$ws.Range("F3").Value = "'0.61"
$ws.Range("H3").Value = "'0.51"
$ws.Range("I3").Value = "'0.54"
$ws.Range("K3").Value = "'0.5"
$ws.Range("L3").Value = "'0.55"
$ws.Range("D4").Value = "'0.29"
$ws.Range("E4").Value = "'0.35"
$ws.Range("I4").Value = "'0.3"
$ws.Range("J4").Value = "'0.29"
$ws.Range("K4").Value = "'0.29"
$ws.Range("E5").Value = "'0.2"
$ws.Range("F5").Value = "'0.1"
$ws.Range("I5").Value = "'0.15"
$ws.Range("J5").Value = "'0.13"
$ws.Range("K5").Value = "'0.17"
$ws.Range("B8").Value = "'0.44"
$ws.Range("C8").Value = "'0.19"
$ws.Range("D8").Value = "'0.32"
$ws.Range("E8").Value = "'0.12"
$ws.Range("F8").Value = "'0.25"
$ws.Range("G8").Value = "'0.11"
$ws.Range("H8").Value = "'0.11"
$ws.Range("I8").Value = "'0.12"
$ws.Range("J8").Value = "'0.06"
$ws.Range("K8").Value = "'0.01"
$ws.Range("L8").Value = "'0.04"
$ws.Range("B9").Value = "'0.31"
$ws.Range("C9").Value = "'0.3"
$ws.Range("D9").Value = "'0.31"
$ws.Range("E9").Value = "'0.26"
$ws.Range("F9").Value = "'0.25"
$ws.Range("G9").Value = "'0.16"
$ws.Range("H9").Value = "'0.13"
$ws.Range("I9").Value = "'0.24"
$ws.Range("J9").Value = "'0.25"
$ws.Range("K9").Value = "'0.12"
$ws.Range("L9").Value = "'0.15"
$ws.Range("B10").Value = "'0.19"
$ws.Range("C10").Value = "'0.38"
$ws.Range("D10").Value = "'0.28"
$ws.Range("E10").Value = "'0.43"
$ws.Range("F10").Value = "'0.37"
$ws.Range("G10").Value = "'0.5"
$ws.Range("H10").Value = "'0.45"
$ws.Range("I10").Value = "'0.48"
$ws.Range("J10").Value = "'0.53"
$ws.Range("K10").Value = "'0.58"
$ws.Range("L10").Value = "'0.55"
$ws.Range("B11").Value = "'0.04"
$ws.Range("C11").Value = "'0.09"
$ws.Range("D11").Value = "'0.08"
$ws.Range("E11").Value = "'0.17"
$ws.Range("F11").Value = "'0.1"
$ws.Range("G11").Value = "'0.2"
$ws.Range("H11").Value = "'0.29"
$ws.Range("I11").Value = "'0.11"
$ws.Range("J11").Value = "'0.14"
$ws.Range("K11").Value = "'0.25"
$ws.Range("L11").Value = "'0.22"
$ws.Range("B13").Value = "'3.55"
$ws.Range("C13").Value = "'2.35"
$ws.Range("D13").Value = "'3.07"
$ws.Range("E13").Value = "'2.48"
$ws.Range("G13").Value = "'2.83"
$ws.Range("H13").Value = "'2.58"
$ws.Range("I13").Value = "'2.39"
$ws.Range("J13").Value = "'2.13"
$ws.Range("K13").Value = "'1.91"
$ws.Range("L13").Value = "'2.25"
$ws.Range("D14").Value = "'0.23"
$ws.Range("G14").Value = "'0.19"
$ws.Range("K14").Value = "'0.3"
$ws.Range("C15").Value = "'0.32"
$ws.Range("D15").Value = "'0.5"
$ws.Range("E15").Value = "'0.3"
$ws.Range("G15").Value = "'0.43"
$ws.Range("H15").Value = "'0.33"
$ws.Range("I15").Value = "'0.29"
$ws.Range("J15").Value = "'0.26"
$ws.Range("K15").Value = "'0.14"
$ws.Range("L15").Value = "'0.25"
$ws.Range("B16").Value = "'0.6"
$ws.Range("C16").Value = "'0.47"
$ws.Range("D16").Value = "'0.44"
$ws.Range("E16").Value = "'0.51"
$ws.Range("F16").Value = "'0.47"
$ws.Range("G16").Value = "'0.55"
$ws.Range("H16").Value = "'0.61"
$ws.Range("I16").Value = "'0.55"
$ws.Range("J16").Value = "'0.38"
$ws.Range("K16").Value = "'0.6"
$ws.Range("L16").Value = "'0.52"
$ws.Range("C18").Value = "'0.27"
$ws.Range("F18").Value = "'0.04"
$ws.Range("A19").Value = "15 - 24"
$ws.Range("B19").Value = "'0.22"
$ws.Range("D19").Value = "'0.16"
$ws.Range("A20").Value = "25 - 49"
$ws.Range("B20").Value = "'0.47"
$ws.Range("E20").Value = "'0.55"
$ws.Range("F20").Value = "'0.55"
$ws.Range("D21").Value = "'0.31"
$ws.Range("E21").Value = "'0.24"
$ws.Range("F21").Value = "'0.25"
$ws.Range("B23").Value = "'0.27"
$ws.Range("C23").Value = "'0.61"
$ws.Range("D23").Value = "'0.17"
$ws.Range("F23").Value = "'0.17"
$ws.Range("G23").Value = "'0.18"
$ws.Range("H23").Value = "'0.33"
$ws.Range("B24").Value = "'0.4"
$ws.Range("C24").Value = "'0.21"
$ws.Range("E24").Value = "'0.4"
$ws.Range("H24").Value = "'0.32"
$ws.Range("E25").Value = "'0.27"
$ws.Range("G25").Value = "'0.31"
$ws.Range("H25").Value = "'0.29"
$ws.Range("B26").Value = "'0.1"
$ws.Range("C26").Value = "'0.07"
$ws.Range("D26").Value = "'0.14"
$ws.Range("E26").Value = "'0.1"
$ws.Range("F26").Value = "'0.1"
$ws.Range("H26").Value = "'0.09"
$ws.Range("B28").Value = "'0.49"
$ws.Range("D28").Value = "'0.58"
$ws.Range("F28").Value = "'0.59"
$ws.Range("G28").Value = "'0.65"
$ws.Range("D29").Value = "'0.62"
$ws.Range("F29").Value = "'0.74"
$ws.Range("G29").Value = "'0.83"
$ws.Range("B30").Value = "'29857"
$ws.Range("C30").Value = "'10973"
$ws.Range("D30").Value = "'3116"
$ws.Range("E30").Value = "'11093"
$ws.Range("F30").Value = "'6702"
$ws.Range("G30").Value = "'9479"
$ws.Range("H30").Value = "'154527"
$ws.Range("I30").Value = "'47087"
$ws.Range("J30").Value = "'141695"
$ws.Range("K30").Value = "'1551834"
$ws.Range("L30").Value = "'34067"
